$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    3 = @{ B = 0.2917716402565462; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 2.591208233317391 }
    4 = @{ B = 0.6606524410359556; C = 1.655778082260271; D = 0.7527432677738641; E = 10.19245300693656; G = 13.26162679800665 }
    5 = @{ B = 1.455362044514542; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 3.754798637575387 }
    6 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    7 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 0.7527432677738641; E = 10.19245300693656; G = 15.88780690183548 }
    8 = @{ B = 0.2917716402565462; C = 0.306821227259698; D = 3.537761648806719; E = 10.19245300693656; G = 14.32880752325952 }
    9 = @{ B = 3.286832544864788; C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
